$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 81, shifting existing rows 81-85 down to 82-86
$ws.Rows.Item(81).Insert()

# Populate the new row 81 with data (copy constant fields from row 82, set changed values)
$ws.Cells.Item(81, 1).Value = 6
$ws.Cells.Item(81, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(81, 3).Value = "Metropolitana"
$ws.Cells.Item(81, 4).Value = 44746
$ws.Cells.Item(81, 5).Value = 13
$ws.Cells.Item(81, 6).Value = 100114007
$ws.Cells.Item(81, 7).Value = "Jengibre"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 100
$ws.Cells.Item(81, 11).Value = 13000
$ws.Cells.Item(81, 12).Value = 15000
$ws.Cells.Item(81, 13).Value = 14100
$ws.Cells.Item(81, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(81, 15).Value = "Perú"
$ws.Cells.Item(81, 16).Value = 1085
$ws.Cells.Item(81, 17).Value = 13
$ws.Cells.Item(81, 18).Value = "Hortaliza"
